# projects_vineet.xlsx edit
#
# The sheet lists projects (col A), location/notes (col B), a "T" marker
# (col D) and, for two rows, a "resource" (col I). This edit:
#   1. Trims the project list down to 6 rows, dropping the "Course5",
#      "Deep Learning / Auto Insurance", "Transactional / vehicle
#      breakdown" and "Tesco Automatic Gap Scan.pdf" rows entirely.
#   2. Drops the now-unused "T" marker column (D) completely.
#   3. Updates the two "HYDERABAD" rows to read "HYDERABAD - data".
#   4. Bolds the whole project-name column, and gives the
#      "Intelligent Clinical Workflow System" row its own (non-wrapping)
#      bold style.
#   5. Tightens rows 1/2 heights and refreshes the active selection/
#      scroll position now that the sheet is shorter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Text fixes (while the original row numbers still apply) --------
$ws.Range("B3").Value2 = "HYDERABAD - data"
$ws.Range("B4").Value2 = "HYDERABAD - data"

# --- 2. Bold the surviving project names --------------------------------
# Rows 1,2,3,8,10 already use Calibri 11 (just not bold) - a plain bold
# toggle reuses that font. Row 9 is still on the old "Segoe UI 10" font,
# so its bold formatting is brought in via a format-only paste from a
# cell that is already Calibri 11 + bold + wrapped, which carries the
# font/alignment pair over in one shot instead of three separate (and
# style-table-bloating) property edits.
$ws.Range("A1:A3").Font.Bold = $true
$ws.Range("A8").Font.Bold = $true
$ws.Range("A10").Font.Bold = $true

$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 4 ("Intelligent Clinical Workflow System") keeps its no-wrap layout
# but becomes bold too.
$ws.Range("A4").Font.Bold = $true

# --- 3. Drop column D (the "T" marker is no longer used) ----------------
$ws.Range("D1:D10").ClearContents()

# --- 4. Row heights for the two tallest remaining rows -------------------
$ws.Rows(1).RowHeight = 30
$ws.Rows(2).RowHeight = 60

# --- 5. Remove the rows that are no longer part of the list --------------
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("9:9").Delete()
$ws.Rows("7:7").Delete()
$ws.Rows("6:6").Delete()
$ws.Rows("5:5").Delete()

# --- 6. Refresh view/selection now that the sheet is only 6 rows tall ----
$ws.Range("B5").Select()
